$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0.901726190476191
$ws.Range("H2").Value = 12

$ws.Range("E3").Value = 0.0969387755102041
$ws.Range("F3").Value = 3.52855737704918
$ws.Range("H3").Value = 42

$ws.Range("E4").Value = 0.320851012088749
$ws.Range("F4").Value = 15.9498800797898
$ws.Range("H4").Value = 84

$ws.Range("E5").Value = 14.9421615952136
$ws.Range("F5").Value = 29.4994372506372
$ws.Range("H5").Value = 57

$ws.Range("F6").Value = 0.364029281616211
$ws.Range("H6").Value = 11

$ws.Range("E7").Value = 0.0419176706827309
$ws.Range("F7").Value = 5.79873585336757
$ws.Range("H7").Value = 25

$ws.Range("E8").Value = 0.611024202184622
$ws.Range("F8").Value = 24.2131707299307
$ws.Range("H8").Value = 37

$ws.Range("E9").Value = 15.6486012328118
$ws.Range("F9").Value = 47.4174548080809
$ws.Range("H9").Value = 16

$ws.Range("F10").Value = 0.236220472440945
$ws.Range("H10").Value = 2

$ws.Range("E11").Value = 0.0426136363636364
$ws.Range("F11").Value = 0.589962121212121
$ws.Range("H11").Value = 5

$ws.Range("E12").Value = 2.35012787723785
$ws.Range("F12").Value = 0.836572890025575
$ws.Range("H12").Value = 19

$ws.Range("E13").Value = 48.3670033670034
$ws.Range("F13").Value = 1.52901023890785
$ws.Range("H13").Value = 19

$ws.Range("E14").Value = 0.0186915887850467
$ws.Range("F14").Value = 0.57314635718891
$ws.Range("H14").Value = 1

$ws.Range("E15").Value = 0.0134347275031686
$ws.Range("F15").Value = 1.22284834704032
$ws.Range("H15").Value = 5

$ws.Range("E16").Value = 0.19156976744186
$ws.Range("F16").Value = 4.78023933880823
$ws.Range("H16").Value = 4

$ws.Range("E17").Value = 4.43903720462544
$ws.Range("F17").Value = 26.2602388746006
$ws.Range("H17").Value = 7
